$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: bold average of column J (the k/n ratio column)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

# Summary rows 14-17
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Format B14 once (bold, size 12, vertical-centered), then propagate the
# exact same style to B15:B17 via copy/paste-format so Excel reuses a
# single cell style instead of minting one per cell.
$ws.Range("B14").Font.Bold = $true
$ws.Range("B14").Font.Size = 12
$ws.Range("B14").VerticalAlignment = -4108

$ws.Range("B14").Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A14:B17").RowHeight = 15.6

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("A14:B17").Select()
